# Update 2p0. Convention change to support multi-axle vehicles
#
# Adds two new vehicle geometry sheets ("Truck_Amandla" and
# "Trailer_Kumanzi") cloned from the existing "Trailer_Thwala" sheet,
# placed either side of it, and adjusts their axle-offset figures.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Create "Truck_Amandla" - a copy of Trailer_Thwala placed directly
#    before it.
# ------------------------------------------------------------------
$thwala = $wb.Worksheets.Item("Trailer_Thwala")
$thwala.Copy($thwala)
$truck = $wb.Worksheets.Item("Trailer_Thwala (2)")
$truck.Name = "Truck_Amandla"

# ------------------------------------------------------------------
# 2. Create "Trailer_Kumanzi" - a copy of Trailer_Thwala placed
#    directly after it.
# ------------------------------------------------------------------
$thwala2 = $wb.Worksheets.Item("Trailer_Thwala")
$thwala2.Copy([System.Reflection.Missing]::Value, $thwala2)
$kumanzi = $wb.Worksheets.Item("Trailer_Thwala (2)")
$kumanzi.Name = "Trailer_Kumanzi"

# ------------------------------------------------------------------
# 3. Both new sheets only need 8 rows (the two spare/template rows
#    9:10 used on Trailer_Thwala are removed).
# ------------------------------------------------------------------
$truck.Rows("9:10").Delete()
$kumanzi.Rows("9:10").Delete()

# ------------------------------------------------------------------
# 4. Re-label the CAD instance name on each new sheet.
# ------------------------------------------------------------------
$truck.Range("H3").Value = "CAD_Truck_Amandla"
$truck.Range("H4").Value = "CAD_Truck_Amandla"

$kumanzi.Range("H3").Value = "CAD_Trailer_Kumanzi"
$kumanzi.Range("H4").Value = "CAD_Trailer_Kumanzi"

# ------------------------------------------------------------------
# 5. Update the sOffset (row 7) / Opacity (row 8) figures for the new
#    multi-axle convention.
# ------------------------------------------------------------------
$truck.Range("F7").Value = 0.6
$truck.Range("G7").Value = 0.8
$truck.Range("H7").Value = 1
$truck.Range("H8").Value = 1

$kumanzi.Range("F7").Value = 1
$kumanzi.Range("G7").Value = 0.75
$kumanzi.Range("H7").Value = 0.055
$kumanzi.Range("H8").Value = 0.5

# ------------------------------------------------------------------
# 6. Fix up selections / active sheet so the view mirrors the
#    original authored workbook (Trailer_Kumanzi ends up active).
#    Trailer_Thwala keeps its original H7 selection untouched.
# ------------------------------------------------------------------
[void]$truck.Activate()
[void]$truck.Range("G23").Select()

[void]$kumanzi.Activate()
[void]$kumanzi.Range("H8").Select()
